$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.736.52'
$ws.Range('E2').Value = '  +0.42%  '
$ws.Range('D3').Value = '1.642.50'
$ws.Range('E3').Value = '  -0.06%  '
$ws.Range('E4').Value = '  +0.15%  '
$ws.Range('D5').Value = '216.46'
$ws.Range('E5').Value = '  +0.40%  '
$ws.Range('D6').Value = '0.501'
$ws.Range('E6').Value = '  -0.97%  '
$ws.Range('E7').Value = '  +0.33%  '
$ws.Range('E8').Value = '  -0.33%  '
$ws.Range('E9').Value = '  +0.30%  '
$ws.Range('D10').Value = '19.11'
$ws.Range('E10').Value = '  -0.50%  '
$ws.Range('E11').Value = '  -0.20%  '
$ws.Range('D12').Value = '1.867.87'
$ws.Range('E12').Value = '  -0.22%  '
$ws.Range('D13').Value = '1.639.27'
$ws.Range('E13').Value = '  -0.35%  '
$ws.Range('E14').Value = '  -1.21%  '
$ws.Range('D15').Value = '0.526'
$ws.Range('E15').Value = '  -0.84%  '
$ws.Range('E16').Value = '  -2.22%  '
$ws.Range('D17').Value = '26.736.94'
$ws.Range('E17').Value = '  +0.31%  '
$ws.Range('D18').Value = '0.0₃0735'
$ws.Range('E18').Value = '  -1.82%  '
$ws.Range('D19').Value = '213.27'
$ws.Range('E19').Value = '  -2.00%  '
$ws.Range('E20').Value = '  +0.19%  '
$ws.Range('E21').Value = '  -0.18%  '
$ws.Range('D22').Value = '2.46'
$ws.Range('E22').Value = '  +13.74%  '
$ws.Range('D23').Value = '6.25'
$ws.Range('E23').Value = '  -0.99%  '
$ws.Range('D24').Value = '9.31'
$ws.Range('E24').Value = '  -2.19%  '
$ws.Range('D25').Value = '145.56'
$ws.Range('E25').Value = '  -0.35%  '
$ws.Range('E26').Value = '  +0.22%  '
$ws.Range('E27').Value = '  -1.70%  '
$ws.Range('E28').Value = '  -0.58%  '
$ws.Range('E29').Value = '  -1.07%  '
$ws.Range('E30').Value = '  -1.67%  '
$ws.Range('E31').Value = '  +0.23%  '
$ws.Range('E32').Value = '  -1.49%  '
$ws.Range('E33').Value = '  -1.95%  '
$ws.Range('D34').Value = '1.291.80'
$ws.Range('E34').Value = '  +1.69%  '
$ws.Range('E35').Value = '  -0.70%  '
$ws.Range('D36').Value = '2.43'
$ws.Range('E36').Value = '  +1.08%  '
$ws.Range('E37').Value = '  -3.26%  '
$ws.Range('E39').Value = '  -1.12%  '
$ws.Range('E40').Value = '  +0.26%  '
$ws.Range('D41').Value = '0.804'
$ws.Range('E41').Value = '  -0.78%  '
$ws.Range('E42').Value = '  -0.54%  '
$ws.Range('B43').Value = 'FraxShare'
$ws.Range('C43').Value = 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs'
$ws.Range('D43').Value = '5.31'
$ws.Range('E43').Value = '  -2.69%  '
$ws.Range('B44').Value = 'RocketPoolETH'
$ws.Range('C44').Value = 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth'
$ws.Range('D44').Value = '1.792.77'
$ws.Range('E44').Value = '  +0.58%  '
$ws.Range('D45').Value = '61.25'
$ws.Range('E45').Value = '  +3.00%  '
$ws.Range('D46').Value = '91.20'
$ws.Range('E46').Value = '  -1.97%  '
$ws.Range('E47').Value = '  +0.27%  '
$ws.Range('B48').Value = 'Cronos'
$ws.Range('C48').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D48').Value = '0.0525'
$ws.Range('E48').Value = '  +1.66%  '
$ws.Range('B49').Value = 'EnergySwap'
$ws.Range('C49').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D49').Value = '7.63'
$ws.Range('E49').Value = '  -1.84%  '
$ws.Range('B50').Value = 'Algorand'
$ws.Range('C50').Value = 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo'
$ws.Range('D50').Value = '0.0973'
$ws.Range('E50').Value = '  +0.02%  '
$ws.Range('B51').Value = 'USDD'
$ws.Range('C51').Value = 'https://coinranking.com/coin/z2PZIKQL7+usdd-usdd'
$ws.Range('D51').Value = '1.01'
$ws.Range('E51').Value = '  +0.45%  '
